# Natmi analysis update following Dr Hou advice:
# - Adds the "ECs" (endothelial cells) sending-cluster group (rows 8-10)
# - Recomputes the specificity-weighted edge statistics for every
#   Sending-cluster x Target-cluster combination (Dhh -> Ptch2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ Row=2; A="ECs"; B="Dhh"; C="Ptch2"; D="ECs"; E=2; F=0.6666666666666666; G=1.853892; H=5.561676; I=0.5711238486747862; J=0.571123848674786; K=3; L=1; M=3.477541333333333; N=10.432624; O=0.2636819480239456; P=0.2636819480239456; Q=6.446986057536001; R=58.02287451782401; S=0.1505950489815008; T=0.1505950489815007 },
  @{ Row=3; A="ECs"; B="Dhh"; C="Ptch2"; D="FAPs"; E=2; F=0.6666666666666666; G=1.853892; H=5.561676; I=0.5711238486747862; J=0.571123848674786; K=3; L=1; M=8.345897; N=25.037691; O=0.6328213435950161; P=0.6328213435950162; Q=15.472391681124; R=139.251525130116; S=0.3614193612775348; T=0.3614193612775348 },
  @{ Row=4; A="ECs"; B="Dhh"; C="Ptch2"; D="sCs"; E=2; F=0.6666666666666666; G=1.853892; H=5.561676; I=0.5711238486747862; J=0.571123848674786; K=3; L=1; M=1.364955333333333; N=4.094866; O=0.1034967083810384; P=0.1034967083810384; Q=2.530479772824; R=22.774317955416; S=0.05910943841575064; T=0.05910943841575062 },
  @{ Row=5; A="FAPs"; B="Dhh"; C="Ptch2"; D="ECs"; E=3; F=1; G=0.4262446666666667; H=1.278734; I=0.1313121230922664; J=0.1313121230922664; K=3; L=1; M=3.477541333333333; N=10.432624; O=0.2636819480239456; P=0.2636819480239456; Q=1.482283446446222; R=13.340551018016; S=0.03462463641612894; T=0.03462463641612894 },
  @{ Row=6; A="FAPs"; B="Dhh"; C="Ptch2"; D="FAPs"; E=3; F=1; G=0.4262446666666667; H=1.278734; I=0.1313121230922664; J=0.1313121230922664; K=3; L=1; M=8.345897; N=25.037691; O=0.6328213435950161; P=0.6328213435950162; Q=3.557394084799334; R=32.01654676319401; S=0.08309711416556217; T=0.08309711416556217 },
  @{ Row=7; A="FAPs"; B="Dhh"; C="Ptch2"; D="sCs"; E=3; F=1; G=0.4262446666666667; H=1.278734; I=0.1313121230922664; J=0.1313121230922664; K=3; L=1; M=1.364955333333333; N=4.094866; O=0.1034967083810384; P=0.1034967083810384; Q=0.5818049310715555; R=5.236244379644; S=0.01359037251057531; T=0.01359037251057531 },
  @{ Row=8; A="sCs"; B="Dhh"; C="Ptch2"; D="ECs"; E=3; F=1; G=0.9659053333333333; H=2.897716; I=0.2975640282329475; J=0.2975640282329475; K=3; L=1; M=3.477541333333333; N=10.432624; O=0.2636819480239456; P=0.2636819480239456; Q=3.358975720753778; R=30.230781486784; S=0.07846226262631595; T=0.07846226262631595 },
  @{ Row=9; A="sCs"; B="Dhh"; C="Ptch2"; D="FAPs"; E=3; F=1; G=0.9659053333333333; H=2.897716; I=0.2975640282329475; J=0.2975640282329475; K=3; L=1; M=8.345897; N=25.037691; O=0.6328213435950161; P=0.6328213435950162; Q=8.061346423750667; R=72.55211781375601; S=0.1883048681519191; T=0.1883048681519192 },
  @{ Row=10; A="sCs"; B="Dhh"; C="Ptch2"; D="sCs"; E=3; F=1; G=0.9659053333333333; H=2.897716; I=0.2975640282329475; J=0.2975640282329475; K=3; L=1; M=1.364955333333333; N=4.094866; O=0.1034967083810384; P=0.1034967083810384; Q=1.318417636228444; R=11.865758726056; S=0.03079689745471243; T=0.03079689745471243 }
)

foreach ($row in $rows) {
    $ws.Range("A" + $row.Row).Value = $row.A
    $ws.Range("B" + $row.Row).Value = $row.B
    $ws.Range("C" + $row.Row).Value = $row.C
    $ws.Range("D" + $row.Row).Value = $row.D
    $ws.Range("E" + $row.Row).Value = $row.E
    $ws.Range("F" + $row.Row).Value = $row.F
    $ws.Range("G" + $row.Row).Value = $row.G
    $ws.Range("H" + $row.Row).Value = $row.H
    $ws.Range("I" + $row.Row).Value = $row.I
    $ws.Range("J" + $row.Row).Value = $row.J
    $ws.Range("K" + $row.Row).Value = $row.K
    $ws.Range("L" + $row.Row).Value = $row.L
    $ws.Range("M" + $row.Row).Value = $row.M
    $ws.Range("N" + $row.Row).Value = $row.N
    $ws.Range("O" + $row.Row).Value = $row.O
    $ws.Range("P" + $row.Row).Value = $row.P
    $ws.Range("Q" + $row.Row).Value = $row.Q
    $ws.Range("R" + $row.Row).Value = $row.R
    $ws.Range("S" + $row.Row).Value = $row.S
    $ws.Range("T" + $row.Row).Value = $row.T
}
